# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) to the left of the existing
# data, shifting the previous Jun_13 / Jun_10 columns two positions to the
# right (B->D, C->E). Bank of America (row 20) receives a new analyst note
# for Jun_15; every other broker stays "UN" (no update) for the two new
# weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand new, blank columns at B:C - this pushes the existing
# Jun_13 column (B) to D and the existing Jun_10 column (C) to E.
$ws.Range("B:C").Insert()

# Header row: new columns get the two newest week labels.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new columns with "UN" (no rating change) for every broker row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Bank of America (row 20) got a new rating note for Jun_15.
$ws.Cells.Item(20, 3).Value = "6/15/2018,Reiterates,Hold,`$163.00"
